# LogTime sprint 7 update: fill in the last three days (30/05, 31/05, 01/06)
# of logged hours on both the "ANLT" and "ANHDT" sheets, then leave the
# selection / active-sheet state the way the author left it.

$wb = $excel.ActiveWorkbook

# --- Sheet "ANLT" (drives Table2) ------------------------------------------
$wsANLT = $wb.Worksheets.Item("ANLT")

# Row 9 ("Meeting") gains three new daily log entries; the SUM() totals row
# (row 13) recalculates automatically from these.
$wsANLT.Range("M9").Value = 4
$wsANLT.Range("N9").Value = 10
$wsANLT.Range("O9").Value = 6

# --- Sheet "ANHDT" (drives Table24) -----------------------------------------
$wsANHDT = $wb.Worksheets.Item("ANHDT")

# Row 2 ("Hien tai admin co the tim kiem order") gains two new daily log
# entries; the SUM() totals row (row 10) recalculates automatically.
$wsANHDT.Range("L2").Value = 2
$wsANHDT.Range("O2").Value = 10

# --- Leave view state as the author left it: ANHDT active/selected with
# O3 selected, ANLT's selection moved to M19.
$wsANLT.Range("M19").Select()

$wsANHDT.Activate()
$wsANHDT.Range("O3").Select()
